$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns A, B, D, E, F, G, H, Q, R between row 13 and row 14
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell13 = $ws.Range("$col`13")
    $cell14 = $ws.Range("$col`14")
    $tmp = $cell13.Value2
    $cell13.Value = $cell14.Value2
    $cell14.Value = $tmp
}
